$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit adds one new weekly record for "Rabanito" at Vega Central Mapocho de
# Santiago (date serial 45009), inserted as the new row 289. All subsequent
# rows (old 289..405) shift down by one row, which Excel's native row-insert
# does for us automatically (values, shared-string refs, row styles, the
# worksheet dimension, etc. all slide down intact).
#
# The brand-new row reuses the same Mercado/Categoria/Calidad metadata as the
# (old) row 289 right below it — only the date changes — so we snapshot those
# values first, insert a blank row above, and then fill it back in.

$origA = $ws.Range("A289").Value2
$origB = $ws.Range("B289").Value2
$origC = $ws.Range("C289").Value2
$origE = $ws.Range("E289").Value2
$origF = $ws.Range("F289").Value2
$origG = $ws.Range("G289").Value2
$origH = $ws.Range("H289").Value2
$origI = $ws.Range("I289").Value2
$origJ = $ws.Range("J289").Value2
$origK = $ws.Range("K289").Value2
$origL = $ws.Range("L289").Value2
$origM = $ws.Range("M289").Value2
$origN = $ws.Range("N289").Value2
$origO = $ws.Range("O289").Value2
$origP = $ws.Range("P289").Value2
$origQ = $ws.Range("Q289").Value2
$origR = $ws.Range("R289").Value2

$ws.Rows("289:289").Insert()

$ws.Range("A289").Value2 = $origA
$ws.Range("B289").Value2 = $origB
$ws.Range("C289").Value2 = $origC
$ws.Range("D289").Value2 = 45009
$ws.Range("E289").Value2 = $origE
$ws.Range("F289").Value2 = $origF
$ws.Range("G289").Value2 = $origG
$ws.Range("H289").Value2 = $origH
$ws.Range("I289").Value2 = $origI
$ws.Range("J289").Value2 = $origJ
$ws.Range("K289").Value2 = $origK
$ws.Range("L289").Value2 = $origL
$ws.Range("M289").Value2 = $origM
$ws.Range("N289").Value2 = $origN
$ws.Range("O289").Value2 = $origO
$ws.Range("P289").Value2 = $origP
$ws.Range("Q289").Value2 = $origQ
$ws.Range("R289").Value2 = $origR
